# Auto commit Tue Jul  5 17:40:01 IST 2022
# - adds a new shared string "4-sigma :"
# - extends the used range down to row 14 with a new "4-sigma :" / 99.99 entry
# - re-expresses the TOTAL row's C11 formula explicitly
# - moves the active selection to O9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TOTAL row C-column: make the formula explicit (SUM(C3:C10)) instead of
# relying on the inherited shared-formula reference.
$ws.Range("C11").Formula = "=SUM(C3:C10)"

# New row 14: a "4-sigma :" label with its value. The column styles
# (col A -> style 2, cols B:K -> style 3) already come from the sheet's
# <cols> defaults, matching rows 3-10 without needing an explicit copy.
$ws.Range("A14").Value = "4-sigma :"
$ws.Range("B14").Value = 99.99

# Move / update the current selection.
$ws.Range("O9").Select()
